# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.055.89'
$ws.Range('E2').Value = '  +5.43%  '
$ws.Range('D3').Value = '1.921.72'
$ws.Range('E3').Value = '  +2.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.82'
$ws.Range('E5').Value = '  +3.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5234'
$ws.Range('E7').Value = '  +3.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4063'
$ws.Range('E8').Value = '  +3.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08472'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.92'
$ws.Range('E10').Value = '  +2.90%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.128'
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.14'
$ws.Range('E12').Value = '  +8.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.370'
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').Value = '1.923.71'
$ws.Range('E14').Value = '  +2.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.380'
$ws.Range('E15').Value = '  +1.69%  '
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '96.22'
$ws.Range('E17').Value = '  +5.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001116'
$ws.Range('E18').Value = '  +1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06741'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.25'
$ws.Range('E20').Value = '  +2.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.065'
$ws.Range('E22').Value = '  +2.18%  '
$ws.Range('D23').Value = '30.072.60'
$ws.Range('E23').Value = '  +5.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.28'
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('E25').Value = '  -1.44%  '
$ws.Range('D26').Value = '2.145.65'
$ws.Range('E26').Value = '  +2.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.19'
$ws.Range('E27').Value = '  +2.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.52'
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.462'
$ws.Range('E29').Value = '  +3.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.86'
$ws.Range('E30').Value = '  +1.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.085'
$ws.Range('E31').Value = '  +4.40%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1062'
$ws.Range('E32').Value = '  +1.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.108'
$ws.Range('E33').Value = '  +5.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.663'
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02521'
$ws.Range('E35').Value = '  +2.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06621'
$ws.Range('E36').Value = '  +1.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2223'
$ws.Range('E37').Value = '  +2.82%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.033'
$ws.Range('E38').Value = '  +1.95%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.237'
$ws.Range('E39').Value = '  +3.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.207'
$ws.Range('E40').Value = '  +2.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6586'
$ws.Range('E41').Value = '  +2.68%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.67'
$ws.Range('E42').Value = '  +5.01%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.247'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6189'
$ws.Range('E44').Value = '  +2.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.21'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.758'
$ws.Range('E46').Value = '  +1.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.076'
$ws.Range('E47').Value = '  +3.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.244'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '125.76'
$ws.Range('E49').Value = '  +3.12%  '
$ws.Range('E50').Value = '  +3.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.63'
$ws.Range('E51').Value = '  +4.15%  '
